# Re-style every table in the deck that currently uses the
# "{B51569CA-6235-427D-B65D-6DDB58766242}" table style so that it uses
# "{84CA5573-C2B2-46E4-BE32-B611FDFAB56D}" instead.
#
# (In this deck that's the tables on slides 14, 15 and 16, but we scan every
# slide/shape so the script isn't dependent on a hard-coded slide list.)

$p = $ppt.ActivePresentation
$oldStyleId = "{B51569CA-6235-427D-B65D-6DDB58766242}"
$newStyleId = "{84CA5573-C2B2-46E4-BE32-B611FDFAB56D}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
